$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '93.629.57'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.36%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.093.62'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -1.66%  '
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '237.53'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -3.61%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '611.44'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.59%  '
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.14'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +2.25%  '
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.387'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.42%  '
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.830'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +10.74%  '
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '3.089.93'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -1.69%  '
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -3.71%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000243'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -4.75%  '
$c.Style = 'Normal'
$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '93.235.49'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.Style = 'Normal'
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '34.72'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -1.86%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.39'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -4.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.665.67'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.64%  '
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.095.33'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -0.97%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '3.63'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -3.21%  '
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '14.76'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -2.37%  '
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.99'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '441.64'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -2.39%  '
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -2.88%  '
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '8.96'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -6.29%  '
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '8.27'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +4.72%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '5.65'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -5.39%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '12.52'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +4.50%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '86.04'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -3.16%  '
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '3.264.02'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -1.22%  '
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -0.25%  '
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.253'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +7.12%  '
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.181'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +6.50%  '
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.126'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -11.70%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '9.14'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -3.46%  '
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +0.78%  '
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '7.93'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -3.22%  '
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.157'
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -11.53%  '
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '25.92'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -2.74%  '
$c.Style = 'Normal'
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.89'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -3.67%  '
$c.Style = 'Normal'
$ws.Range('B40').Value = 'MantraDAO'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.87'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -9.55%  '
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.448'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +0.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '23.99'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +7.94%  '
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '473.22'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -4.72%  '
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -3.98%  '
$c.Style = 'Normal'
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.20'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -8.11%  '
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '159.98'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.696'
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.57%  '
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.84'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -5.40%  '
$c.Style = 'Normal'
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '43.79'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.93%  '
$c.Style = 'Normal'
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.31'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -5.62%  '
$c.Style = 'Normal'
